$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new row 33 averages (5 * rubric weight) under columns E/F/G
$ws.Range("E33").Formula = "=5*0.89"
$ws.Range("F33").Formula = "=5*0.915"
$ws.Range("G33").Formula = "=5*0.93"

# Scroll the view so B13 is the top-left visible cell and C4 is selected
$win = $excel.ActiveWindow
$win.ScrollRow = 13
$win.ScrollColumn = 2
$ws.Range("C4").Select()
